$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 25000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 278.4
$ws.Range("I31").Value = 278.4
$ws.Range("K31").Value = 835.1999999999999
$ws.Range("M31").Value = -605.1999999999999
$ws.Range("H33").Value = 263.7857
$ws.Range("I33").Value = 263.7857
$ws.Range("K33").Value = 263.7857
$ws.Range("M33").Value = -34.78570000000002
$ws.Range("H39").Value = 321.6875
$ws.Range("I39").Value = 213
$ws.Range("K39").Value = 639
$ws.Range("M39").Value = -343
$ws.Range("H61").Value = 4751.263
$ws.Range("I61").Value = 4751.263
$ws.Range("K61").Value = 14253.789
$ws.Range("M61").Value = -14081.789
$ws.Range("H62").Value = 9811655
$ws.Range("I62").Value = 18525794
$ws.Range("K62").Value = 18525794
$ws.Range("M62").Value = -18525170
$ws.Range("H65").Value = 9811655
$ws.Range("I65").Value = 18525794
$ws.Range("K65").Value = 92628970
$ws.Range("M65").Value = -92625850
$ws.Range("H70").Value = 14214.538
$ws.Range("J70").Value = 19499.889
$ws.Range("L70").Value = 58499.667
$ws.Range("N70").Value = -59039.667
$ws.Range("H73").Value = 14214.538
$ws.Range("J73").Value = 19499.889
$ws.Range("L73").Value = 58499.667
$ws.Range("N73").Value = -60371.667
$ws.Range("H86").Value = 6774.95
$ws.Range("I86").Value = 7156.5713
$ws.Range("J86").Value = 6569.4614
$ws.Range("K86").Value = 7156.5713
$ws.Range("L86").Value = 6569.4614
$ws.Range("M86").Value = -6033.5713
$ws.Range("N86").Value = -8815.4614
$ws.Range("H89").Value = 6774.95
$ws.Range("I89").Value = 7156.5713
$ws.Range("J89").Value = 6569.4614
$ws.Range("K89").Value = 35782.85649999999
$ws.Range("L89").Value = 32847.307
$ws.Range("M89").Value = -30166.85649999999
$ws.Range("N89").Value = -44079.307
$ws.Range("H98").Value = 1256075.1
$ws.Range("I98").Value = 800
$ws.Range("K98").Value = 800
$ws.Range("M98").Value = 698
$ws.Range("H106").Value = 3206.0715
$ws.Range("I106").Value = 3144.5386
$ws.Range("K106").Value = 3144.5386
$ws.Range("M106").Value = -2513.5386
$ws.Range("H107").Value = 1630.6072
$ws.Range("I107").Value = 1833.1
$ws.Range("K107").Value = 1833.1
$ws.Range("M107").Value = 86.90000000000009
$ws.Range("H122").Value = 1256075.1
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50
$ws.Range("H127").Value = 77736470
$ws.Range("I127").Value = 1010.7143
$ws.Range("K127").Value = 3032.1429
$ws.Range("M127").Value = 1927.8571
$ws.Range("H132").Value = 3081.8572
$ws.Range("I132").Value = 3126.1904
$ws.Range("K132").Value = 9378.5712
$ws.Range("M132").Value = -6848.5712
$ws.Range("H137").Value = 3162372
$ws.Range("I137").Value = 5618991.5
$ws.Range("K137").Value = 16856974.5
$ws.Range("M137").Value = -16854424.5
$ws.Range("H138").Value = 3694.9697
$ws.Range("I138").Value = 1911.1666
$ws.Range("J138").Value = 4714.2856
$ws.Range("K138").Value = 5733.4998
$ws.Range("L138").Value = 14142.8568
$ws.Range("M138").Value = -593.4997999999996
$ws.Range("N138").Value = -24422.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4943.364
$ws.Range("I45").Value = 4943.1
$ws.Range("K45").Value = 4943.1
$ws.Range("M45").Value = -4566.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 40016
$ws.Range("I24").Value = 40016
$ws.Range("K24").Value = 40016
$ws.Range("M24").Value = -39781
$ws.Range("H25").Value = 511.33334
$ws.Range("I25").Value = 511.33334
$ws.Range("K25").Value = 511.33334
$ws.Range("M25").Value = -276.33334
$ws.Range("H29").Value = 18
$ws.Range("J29").Value = 18
$ws.Range("L29").Value = 18
$ws.Range("N29").Value = -596
$ws.Range("H134").Value = 6035.4614
$ws.Range("I134").Value = 6035.4614
$ws.Range("K134").Value = 18106.3842
$ws.Range("M134").Value = -15571.3842

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 266246.94
$ws.Range("I31").Value = 336446.4
$ws.Range("J31").Value = 2999
$ws.Range("K31").Value = 336446.4
$ws.Range("L31").Value = 2999
$ws.Range("M31").Value = -336151.4
$ws.Range("N31").Value = -3589
$ws.Range("H34").Value = 266246.94
$ws.Range("I34").Value = 336446.4
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 336446.4
$ws.Range("L34").Value = 2999
$ws.Range("M34").Value = -336244.4
$ws.Range("N34").Value = -3403
$ws.Range("H86").Value = 1006056.4
$ws.Range("I86").Value = 2005800.6
$ws.Range("J86").Value = 6312.2
$ws.Range("K86").Value = 2005800.6
$ws.Range("L86").Value = 6312.2
$ws.Range("M86").Value = -2004677.6
$ws.Range("N86").Value = -8558.200000000001
$ws.Range("H89").Value = 1006056.4
$ws.Range("I89").Value = 2005800.6
$ws.Range("J89").Value = 6312.2
$ws.Range("K89").Value = 10029003
$ws.Range("L89").Value = 31561
$ws.Range("M89").Value = -10023387
$ws.Range("N89").Value = -42793
$ws.Range("H105").Value = 5304.184
$ws.Range("I105").Value = 2098.0625
$ws.Range("K105").Value = 2098.0625
$ws.Range("M105").Value = -351.0625
$ws.Range("H132").Value = 6649.915
$ws.Range("I132").Value = 6170.8203
$ws.Range("J132").Value = 8985.5
$ws.Range("K132").Value = 18512.4609
$ws.Range("L132").Value = 26956.5
$ws.Range("M132").Value = -15982.4609
$ws.Range("N132").Value = -32016.5
$ws.Range("H134").Value = 2907.9736
$ws.Range("I134").Value = 2652.861
$ws.Range("K134").Value = 7958.583
$ws.Range("M134").Value = -5423.583

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 250411.53
$ws.Range("I4").Value = 150441.22
$ws.Range("K4").Value = 451323.66
$ws.Range("M4").Value = -451211.66
$ws.Range("H137").Value = 83334090
$ws.Range("I137").Value = 1006.6667
$ws.Range("K137").Value = 3020.0001
$ws.Range("M137").Value = 2079.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3632.85
$ws.Range("I80").Value = 3427.8235
$ws.Range("J80").Value = 3784.3914
$ws.Range("K80").Value = 3427.8235
$ws.Range("L80").Value = 3784.3914
$ws.Range("M80").Value = -2429.8235
$ws.Range("N80").Value = -5780.3914
$ws.Range("H83").Value = 3632.85
$ws.Range("I83").Value = 3427.8235
$ws.Range("J83").Value = 3784.3914
$ws.Range("K83").Value = 17139.1175
$ws.Range("L83").Value = 18921.957
$ws.Range("M83").Value = -12147.1175
$ws.Range("N83").Value = -28905.957

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 18999.75
$ws.Range("I4").Value = 18999
$ws.Range("K4").Value = 18999
$ws.Range("M4").Value = -18886
$ws.Range("H22").Value = 1322
$ws.Range("J22").Value = 1442.7142
$ws.Range("L22").Value = 1442.7142
$ws.Range("N22").Value = -2032.7142
$ws.Range("H27").Value = 1322
$ws.Range("J27").Value = 1442.7142
$ws.Range("L27").Value = 1442.7142
$ws.Range("N27").Value = -1656.7142
$ws.Range("H28").Value = 18999.75
$ws.Range("I28").Value = 18999
$ws.Range("K28").Value = 18999
$ws.Range("M28").Value = -18767
$ws.Range("H37").Value = 18999.75
$ws.Range("I37").Value = 18999
$ws.Range("K37").Value = 18999
$ws.Range("M37").Value = -18892
$ws.Range("H132").Value = 6903.24
$ws.Range("I132").Value = 4729.2354
$ws.Range("K132").Value = 14187.7062
$ws.Range("M132").Value = -11657.7062
$ws.Range("H136").Value = 3854.7827
$ws.Range("I136").Value = 3858.611
$ws.Range("K136").Value = 11575.833
$ws.Range("M136").Value = -9025.832999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9575
$ws.Range("J62").Value = 9832.5
$ws.Range("L62").Value = 9832.5
$ws.Range("N62").Value = -11080.5
$ws.Range("H65").Value = 9575
$ws.Range("J65").Value = 9832.5
$ws.Range("L65").Value = 49162.5
$ws.Range("N65").Value = -55402.5
$ws.Range("H81").Value = 3553.25
$ws.Range("I81").Value = 2111.7693
$ws.Range("J81").Value = 9799.666999999999
$ws.Range("K81").Value = 4223.5386
$ws.Range("L81").Value = 19599.334
$ws.Range("M81").Value = -3162.5386
$ws.Range("N81").Value = -21721.334
$ws.Range("H84").Value = 3553.25
$ws.Range("I84").Value = 2111.7693
$ws.Range("J84").Value = 9799.666999999999
$ws.Range("K84").Value = 21117.693
$ws.Range("L84").Value = 97996.67
$ws.Range("M84").Value = -15813.693
$ws.Range("N84").Value = -108604.67
